$d = $word.ActiveDocument

# 1. Change the title paragraph's style from Heading1 to Title
$d.Paragraphs(1).Style = "Title"

# 2. Insert a new paragraph with "Initial results:" text before the
#    bookmark paragraph (which holds the lone _GoBack bookmark).
$bookmarkPara = $d.Paragraphs(3)
$bookmarkPara.Range.InsertParagraphBefore()
$d.Paragraphs(3).Range.Text = "Initial results:"

# 3. Switch the page to landscape orientation.
$d.PageSetup.Orientation = 1

# 4. Applying the built-in "Title" style above auto-mints a bare-bones
#    style entry. Flesh it out (and add its linked "Title Char" style)
#    to match Word's real built-in Title definition.
$titleStyle = $d.Styles("Title")
$titleStyle.NextParagraphStyle = "Normal"
$titleStyle.Priority = 10
$titleStyle.QuickStyle = $true
$titleStyle.ParagraphFormat.SpaceAfter = 0
$titleStyle.ParagraphFormat.LineSpacingRule = 0
$titleStyle.NoSpaceBetweenParagraphsOfSameStyle = $true
$titleStyle.Font.Size = 28
$titleStyle.Font.SizeBi = 28
$titleStyle.Font.Kerning = 14
$titleStyle.Font.Spacing = -0.5

$titleCharStyle = $d.Styles.Add("TitleChar", 2)
$titleCharStyle.NameLocal = "Title Char"
$titleCharStyle.BaseStyle = "DefaultParagraphFont"
$titleCharStyle.Priority = 10
$titleCharStyle.Font.Size = 28
$titleCharStyle.Font.SizeBi = 28
$titleCharStyle.Font.Kerning = 14
$titleCharStyle.Font.Spacing = -0.5

$titleStyle.LinkStyle = "TitleChar"
$titleCharStyle.LinkStyle = "Title"

Write-Output "done"
